# Coastal Surface Piercing Profilers: re-point the CP05MOAS-GL001 mooring
# cal-info sheet at the GL335 mooring (corrected instrument reference
# designators), update deployment numbers, and move the saved selections.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# "Moorings" sheet
# ----------------------------------------------------------------------
$moorings = $wb.Worksheets.Item("Moorings")

# Ref Des column: CP05MOAS-GL001 -> CP05MOAS-GL335
$moorings.Range("A2").Value = "CP05MOAS-GL335"

# Deployment Number: 2 -> 1
$moorings.Range("C2").Value = 1

# ----------------------------------------------------------------------
# "Asset_Cal_Info" sheet
# ----------------------------------------------------------------------
$cal = $wb.Worksheets.Item("Asset_Cal_Info")

# Corrected instrument reference designators: GL001 -> GL335,
# and deployment number 2 -> 1 for every data row.
$cal.Range("A2").Value = "CP05MOAS-GL335-01-ADCPAM000"
$cal.Range("A3").Value = "CP05MOAS-GL335-01-ADCPAM000"
$cal.Range("A4").Value = "CP05MOAS-GL335-01-ADCPAM000"
$cal.Range("A5").Value = "CP05MOAS-GL335-01-ADCPAM000"
$cal.Range("C2").Value = 1
$cal.Range("C3").Value = 1
$cal.Range("C4").Value = 1
$cal.Range("C5").Value = 1

$cal.Range("A7").Value = "CP05MOAS-GL335-02-FLORTM000"
$cal.Range("A8").Value = "CP05MOAS-GL335-02-FLORTM000"
$cal.Range("A9").Value = "CP05MOAS-GL335-02-FLORTM000"
$cal.Range("A10").Value = "CP05MOAS-GL335-02-FLORTM000"
$cal.Range("C7").Value = 1
$cal.Range("C8").Value = 1
$cal.Range("C9").Value = 1
$cal.Range("C10").Value = 1

$cal.Range("A12").Value = "CP05MOAS-GL335-03-CTDGVM000"
$cal.Range("C12").Value = 1

$cal.Range("A14").Value = "CP05MOAS-GL335-04-DOSTAM000"
$cal.Range("C14").Value = 1

$cal.Range("A16").Value = "CP05MOAS-GL335-05-PARADM000"
$cal.Range("C16").Value = 1

$cal.Range("A18").Value = "CP05MOAS-GL335-00-ENG000000"
$cal.Range("C18").Value = 1

# Saved selection moves from F18 to C18
[void]$cal.Range("C18").Select()

# "Moorings" remains the active sheet/tab; its own saved selection moves
# from D2 to B12 (re-activate it last so it stays the selected tab).
$moorings.Activate()
[void]$moorings.Range("B12").Select()
